$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    "G2" = 37.84230566666667
    "H2" = 113.526917
    "I2" = 0.5048163536019187
    "J2" = 0.5048163536019187
    "M2" = 0.6795853333333334
    "N2" = 2.038756
    "O2" = 0.08284139605799233
    "P2" = 0.08284139605799234
    "Q2" = 25.71707591058356
    "R2" = 231.453683195252
    "S2" = 0.04181969148528805
    "T2" = 0.04181969148528806
    "G3" = 37.84230566666667
    "H3" = 113.526917
    "I3" = 0.5048163536019187
    "J3" = 0.5048163536019187
    "O3" = 0.3136748993401273
    "P3" = 0.3136748993401273
    "Q3" = 97.37645164656112
    "R3" = 876.38806481905
    "S3" = 0.1583482189013319
    "T3" = 0.1583482189013319
    "G4" = 37.84230566666667
    "H4" = 113.526917
    "I4" = 0.5048163536019187
    "J4" = 0.5048163536019187
    "M4" = 2.840162333333334
    "N4" = 8.520487000000001
    "O4" = 0.3462155540800247
    "P4" = 0.3462155540800247
    "Q4" = 107.4782911609532
    "R4" = 967.304620448579
    "S4" = 0.174775273570946
    "T4" = 0.174775273570946
    "G5" = 37.84230566666667
    "H5" = 113.526917
    "I5" = 0.5048163536019187
    "J5" = 0.5048163536019187
    "M5" = 2.110486666666667
    "N5" = 6.33146
    "O5" = 0.2572681505218555
    "P5" = 0.2572681505218555
    "Q5" = 79.86568154542445
    "R5" = 718.7911339088199
    "S5" = 0.1298731696443527
    "T5" = 0.1298731696443527
    "G6" = 15.88630666666666
    "H6" = 47.65891999999999
    "I6" = 0.2119233292577262
    "J6" = 0.2119233292577262
    "M6" = 0.6795853333333334
    "N6" = 2.038756
    "O6" = 0.08284139605799233
    "P6" = 0.08284139605799234
    "Q6" = 10.79610101150222
    "R6" = 97.16490910352
    "S6" = 0.01755602445296761
    "T6" = 0.01755602445296761
    "G7" = 15.88630666666666
    "H7" = 47.65891999999999
    "I7" = 0.2119233292577262
    "J7" = 0.2119233292577262
    "O7" = 0.3136748993401273
    "P7" = 0.3136748993401273
    "Q7" = 40.87890908644444
    "R7" = 367.9101817779999
    "S7" = 0.06647502897274191
    "T7" = 0.06647502897274191
    "G8" = 15.88630666666666
    "H8" = 47.65891999999999
    "I8" = 0.2119233292577262
    "J8" = 0.2119233292577262
    "M8" = 2.840162333333334
    "N8" = 8.520487000000001
    "O8" = 0.3462155540800247
    "P8" = 0.3462155540800247
    "Q8" = 45.11968981044889
    "R8" = 406.07720829404
    "S8" = 0.07337115286144719
    "T8" = 0.0733711528614472
    "G9" = 15.88630666666666
    "H9" = 47.65891999999999
    "I9" = 0.2119233292577262
    "J9" = 0.2119233292577262
    "M9" = 2.110486666666667
    "N9" = 6.33146
    "O9" = 0.2572681505218555
    "P9" = 0.2572681505218555
    "Q9" = 33.52783840257777
    "R9" = 301.7505456232
    "S9" = 0.05452112297056945
    "T9" = 0.05452112297056946
    "G10" = 18.76675533333333
    "H10" = 56.300266
    "I10" = 0.2503485141672444
    "J10" = 0.2503485141672445
    "M10" = 0.6795853333333334
    "N10" = 2.038756
    "O10" = 0.08284139605799233
    "P10" = 0.08284139605799234
    "Q10" = 12.75361167878845
    "R10" = 114.782505109096
    "S10" = 0.0207392204146586
    "T10" = 0.02073922041465861
    "G11" = 18.76675533333333
    "H11" = 56.300266
    "I11" = 0.2503485141672444
    "J11" = 0.2503485141672445
    "O11" = 0.3136748993401273
    "P11" = 0.3136748993401273
    "Q11" = 48.29092760298889
    "R11" = 434.6183484269
    "S11" = 0.07852804498136082
    "T11" = 0.07852804498136083
    "G12" = 18.76675533333333
    "H12" = 56.300266
    "I12" = 0.2503485141672444
    "J12" = 0.2503485141672445
    "M12" = 2.840162333333334
    "N12" = 8.520487000000001
    "O12" = 0.3462155540800247
    "P12" = 0.3462155540800247
    "Q12" = 53.30063161661578
    "R12" = 479.7056845495421
    "S12" = 0.08667454954552345
    "T12" = 0.08667454954552346
    "G13" = 18.76675533333333
    "H13" = 56.300266
    "I13" = 0.2503485141672444
    "J13" = 0.2503485141672445
    "M13" = 2.110486666666667
    "N13" = 6.33146
    "O13" = 0.2572681505218555
    "P13" = 0.2572681505218555
    "Q13" = 39.60698690759556
    "R13" = 356.46288216836
    "S13" = 0.06440669922570151
    "T13" = 0.06440669922570152
    "G14" = 2.467151666666667
    "H14" = 7.401455
    "I14" = 0.03291180297311068
    "J14" = 0.03291180297311068
    "M14" = 0.6795853333333334
    "N14" = 2.038756
    "O14" = 0.08284139605799233
    "P14" = 0.08284139605799234
    "Q14" = 1.676640087775556
    "R14" = 15.08976078998
    "S14" = 0.002726459705078071
    "T14" = 0.002726459705078072
    "G15" = 2.467151666666667
    "H15" = 7.401455
    "I15" = 0.03291180297311068
    "J15" = 0.03291180297311068
    "O15" = 0.3136748993401273
    "P15" = 0.3136748993401273
    "Q15" = 6.348515787861111
    "R15" = 57.13664209075
    "S15" = 0.01032360648469259
    "T15" = 0.01032360648469259
    "G16" = 2.467151666666667
    "H16" = 7.401455
    "I16" = 0.03291180297311068
    "J16" = 0.03291180297311068
    "M16" = 2.840162333333334
    "N16" = 8.520487000000001
    "O16" = 0.3462155540800247
    "P16" = 0.3462155540800247
    "Q16" = 7.007111234287223
    "R16" = 63.06400110858501
    "S16" = 0.01139457810210812
    "T16" = 0.01139457810210812
    "G17" = 2.467151666666667
    "H17" = 7.401455
    "I17" = 0.03291180297311068
    "J17" = 0.03291180297311068
    "M17" = 2.110486666666667
    "N17" = 6.33146
    "O17" = 0.2572681505218555
    "P17" = 0.2572681505218555
    "Q17" = 5.206890697144444
    "R17" = 46.8620162743
    "S17" = 0.008467158681231891
    "T17" = 0.008467158681231891
}

foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}
